$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-24 21:01:37"
$wsZhCn.Range("H3").Value = "2016-03-24 21:02:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-24 21:01:45"
$wsDeDe.Range("H3").Value = "2016-03-24 21:02:20"
